# Mark the "CC validator" project (row 7) as finished by recording its
# finish date in column F ("Date Finished"), matching the formatting
# already used for the other Date Finished / Date Started cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List of Projects")

# Copy the date formatting from the neighboring "Date Started" cell (E7)
# so the new cell picks up the existing date number format/style instead
# of creating a brand new style entry.
$ws.Range("E7").Copy() | Out-Null
$ws.Range("F7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Record the finish date (2019-02-23) for the CC validator project.
$ws.Range("F7").Value = 43519

# Reflect where the user's cursor ended up after finishing the edit.
$ws.Range("I9").Select() | Out-Null
